$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two THT electrolytic capacitor rows (100u and 1000u 16V) that are
# no longer present in the exported manufacturing BOM. Deleting the rows shifts
# everything below up by two and the shared formulas / SUM range adjust
# automatically.
$ws.Range("A8:A9").EntireRow.Delete()

# Restore the view/selection state captured when the data was exported.
$ws.Range("M40").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
